$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F/G column values (AgTests / AgPosit corrections)
$ws.Range("F483").Value = 66725
$ws.Range("F513").Value = 11241
$ws.Range("F530").Value = 13296
$ws.Range("F544").Value = 14787
$ws.Range("F547").Value = 14134
$ws.Range("F548").Value = 17350
$ws.Range("F551").Value = 18226
$ws.Range("F555").Value = 21739
$ws.Range("F558").Value = 24806
$ws.Range("F559").Value = 22526
$ws.Range("F561").Value = 24509
$ws.Range("F562").Value = 27118
$ws.Range("F565").Value = 29162
$ws.Range("F568").Value = 24101
$ws.Range("F572").Value = 33539
$ws.Range("F573").Value = 27014
$ws.Range("F575").Value = 26385
$ws.Range("F576").Value = 28998
$ws.Range("F579").Value = 32834
$ws.Range("F582").Value = 26130
$ws.Range("F583").Value = 29429
$ws.Range("F586").Value = 33856
$ws.Range("F589").Value = 26056
$ws.Range("F590").Value = 29017
$ws.Range("F593").Value = 37235
$ws.Range("F596").Value = 29371
$ws.Range("F597").Value = 29578
$ws.Range("F600").Value = 40157
$ws.Range("F602").Value = 30085
$ws.Range("F603").Value = 32022
$ws.Range("F604").Value = 29927
$ws.Range("F606").Value = 14315
$ws.Range("F607").Value = 10822
$ws.Range("F608").Value = 45973
$ws.Range("F609").Value = 36373
$ws.Range("F610").Value = 33881
$ws.Range("F611").Value = 34071
$ws.Range("G611").Value = 2125
$ws.Range("F613").Value = 21595
$ws.Range("F614").Value = 47614
$ws.Range("F615").Value = 36622
$ws.Range("G615").Value = 2347
$ws.Range("F616").Value = 37597
$ws.Range("F617").Value = 38789
$ws.Range("G617").Value = 2591
$ws.Range("F618").Value = 37573
$ws.Range("G618").Value = 2653
$ws.Range("F619").Value = 17676
$ws.Range("G619").Value = 1881
$ws.Range("F620").Value = 25471
$ws.Range("G620").Value = 2366
$ws.Range("F621").Value = 55404
$ws.Range("G621").Value = 4070
$ws.Range("F622").Value = 40712
$ws.Range("G622").Value = 2968
$ws.Range("F623").Value = 14798
$ws.Range("G623").Value = 1543

# Append new daily rows for 2021-11-18 through 2021-11-21
$ws.Range("A624").Value = 44518
$ws.Range("B624").Value = 593242
$ws.Range("C624").Value = 23440
$ws.Range("D624").Value = 7418
$ws.Range("E624").Value = 13781
$ws.Range("F624").Value = 50246
$ws.Range("G624").Value = 3918

$ws.Range("A625").Value = 44519
$ws.Range("B625").Value = 602413
$ws.Range("C625").Value = 27267
$ws.Range("D625").Value = 9171
$ws.Range("E625").Value = 13818
$ws.Range("F625").Value = 38338
$ws.Range("G625").Value = 3056

$ws.Range("A626").Value = 44520
$ws.Range("B626").Value = 610140
$ws.Range("C626").Value = 22096
$ws.Range("D626").Value = 7727
$ws.Range("E626").Value = 13861
$ws.Range("F626").Value = 17416
$ws.Range("G626").Value = 1798

$ws.Range("A627").Value = 44521
$ws.Range("B627").Value = 614684
$ws.Range("C627").Value = 13025
$ws.Range("D627").Value = 4544
$ws.Range("E627").Value = 13919
$ws.Range("F627").Value = 25712
$ws.Range("G627").Value = 2101

